# Weekly update: a new price-survey row for "Ají" (Agrícola del Norte S.A.
# de Arica) is inserted at row 54, pushing the previously-existing rows
# 54..77 down to 55..78 (dimension grows from A1:R77 to A1:R78).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 54; this shifts rows
# 54..77 down to 55..78 automatically.
$ws.Rows.Item(54).Insert()

# Populate the new row 54 with this week's data point (same fixed
# metadata columns as every other row in this sheet: Mercado, Región,
# Codreg, Categoría ID/Categoría, Unidad de comercialización, Origen,
# Kg o Unidades, Clasificación).
$ws.Cells.Item(54, 1).Value = 1
$ws.Cells.Item(54, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(54, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(54, 4).Value = 44704
$ws.Cells.Item(54, 5).Value = 15
$ws.Cells.Item(54, 6).Value = 100112021
$ws.Cells.Item(54, 7).Value = "Ají"
$ws.Cells.Item(54, 8).Value = "Inferno"
$ws.Cells.Item(54, 9).Value = "Primera"
$ws.Cells.Item(54, 10).Value = 130
$ws.Cells.Item(54, 11).Value = 19000
$ws.Cells.Item(54, 12).Value = 20000
$ws.Cells.Item(54, 13).Value = 19500
$ws.Cells.Item(54, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(54, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(54, 16).Value = 1300
$ws.Cells.Item(54, 17).Value = 15
$ws.Cells.Item(54, 18).Value = "Hortaliza"

# Keep the date column's display format consistent with the rest of
# column D (custom "YYYY-MM-DD HH:MM:SS" format already used below).
$ws.Cells.Item(54, 4).NumberFormat = $ws.Cells.Item(55, 4).NumberFormat
